# Apply updated crypto price/volume data per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.958.81"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").Value = "2.469.46"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'518.41"
$ws.Range("E5").Value = "  -3.54%  "
$ws.Range("D6").Value = "'130.99"
$ws.Range("E6").Value = "  -4.43%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  -2.56%  "
$ws.Range("D9").Value = "'0.0991"
$ws.Range("E9").Value = "  -2.12%  "
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Value = "'0.343"
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("D13").Value = "2.908.87"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("D14").Value = "57.877.84"
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("D15").Value = "'22.27"
$ws.Range("E15").Value = "  -3.16%  "
$ws.Range("E16").Value = "  -2.20%  "
$ws.Range("D17").Value = "2.467.45"
$ws.Range("E17").Value = "  -1.70%  "
$ws.Range("D18").Value = "'10.82"
$ws.Range("E18").Value = "  -2.65%  "
$ws.Range("D19").Value = "'4.17"
$ws.Range("E19").Value = "  -2.25%  "
$ws.Range("D20").Value = "'318.71"
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "'5.74"
$ws.Range("E22").Value = "  -3.46%  "
$ws.Range("D23").Value = "'64.05"
$ws.Range("E23").Value = "  -2.79%  "
$ws.Range("E24").Value = "  -2.93%  "
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").Value = "'0.159"
$ws.Range("E26").Value = "  -3.22%  "
$ws.Range("D27").Value = "'7.33"
$ws.Range("E27").Value = "  -2.82%  "
$ws.Range("D28").Value = "0.0₃0753"
$ws.Range("E28").Value = "  -2.40%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'1.69"
$ws.Range("E29").Value = "  -4.23%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'166.16"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("E31").Value = "  -5.76%  "
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("E35").Value = "  -2.07%  "
$ws.Range("E36").Value = "  -9.60%  "
$ws.Range("D37").Value = "'3.98"
$ws.Range("E37").Value = "  -3.08%  "
$ws.Range("E38").Value = "  -4.48%  "
$ws.Range("D39").Value = "'0.790"
$ws.Range("E39").Value = "  -2.48%  "
$ws.Range("E40").Value = "  -4.13%  "
$ws.Range("D41").Value = "'273.33"
$ws.Range("E41").Value = "  -3.86%  "
$ws.Range("D42").Value = "'4.99"
$ws.Range("E42").Value = "  -4.48%  "
$ws.Range("E43").Value = "  -2.48%  "
$ws.Range("D44").Value = "'126.42"
$ws.Range("E44").Value = "  -4.79%  "
$ws.Range("D45").Value = "'0.0905"
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("E46").Value = "  -3.75%  "
$ws.Range("E47").Value = "  -3.11%  "
$ws.Range("D48").Value = "'17.10"
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("D49").Value = "1.734.38"
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("E51").Value = "  -0.96%  "
